$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report row was inserted before the existing row 34,
# shifting all subsequent rows (old 34-45) down to (35-46).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record.
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 44449
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100107
$ws.Range("H34").Value = "Otros"
$ws.Range("I34").Value = 100107002
$ws.Range("J34").Value = "Chirimoya"
$ws.Range("K34").Value = "Cultivar IV Región"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 50
$ws.Range("N34").Value = 3000
$ws.Range("O34").Value = 3000
$ws.Range("P34").Value = 3000
$ws.Range("Q34").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R34").Value = "Provincia del Elquí"
$ws.Range("S34").Value = 3000
$ws.Range("T34").Value = 1

# Make sure the date cell keeps the date number format used by the rest of
# the column (style index 2 in the original workbook).
$ws.Range("D34").NumberFormat = $ws.Range("D35").NumberFormat
